# Update cryptocurrency Price (column D) and Volume(1h) (column E) values
# for rows 2-51, per the "Updated cryptos list" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'28.064.76"
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -1.22%  '
$ws.Range('D3').Value = "'1.794.42"
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +0.15%  '
$ws.Range('D4').Value = "'1.002"
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'317.58"
$ws.Range('D5').Style = "Normal"
$ws.Range('E6').Value = '  +0.11%  '
$ws.Range('D7').Value = "'0.5384"
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -1.67%  '
$ws.Range('E8').Value = '  -1.20%  '
$ws.Range('D9').Value = "'0.07446"
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -1.64%  '
$ws.Range('E10').Value = '  -1.71%  '
$ws.Range('D11').Value = "'1.093"
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -2.45%  '
$ws.Range('E12').Value = '  +0.06%  '
$ws.Range('D13').Value = "'20.52"
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -2.65%  '
$ws.Range('E14').Value = '  -1.17%  '
$ws.Range('D15').Value = "'1.799.20"
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.29%  '
$ws.Range('D16').Value = "'7.233"
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -2.10%  '
$ws.Range('E17').Value = '  -2.53%  '
$ws.Range('D18').Value = "'0.00001060"
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.72%  '
$ws.Range('D19').Value = "'0.06481"
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.32%  '
$ws.Range('E20').Value = '  -0.04%  '
$ws.Range('D21').Value = "'17.24"
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.62%  '
$ws.Range('D22').Value = "'5.895"
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.95%  '
$ws.Range('D23').Value = "'28.090.18"
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -1.05%  '
$ws.Range('E24').Value = '  -2.11%  '
$ws.Range('D25').Value = "'2.089"
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -1.49%  '
$ws.Range('D26').Value = "'155.37"
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -2.81%  '
$ws.Range('E27').Value = '  -2.13%  '
$ws.Range('D28').Value = "'2.001.68"
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('E29').Value = '  -5.04%  '
$ws.Range('D30').Value = "'121.11"
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.48%  '
$ws.Range('D31').Value = "'1.116"
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.44%  '
$ws.Range('D32').Value = "'0.1061"
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +3.40%  '
$ws.Range('D33').Value = "'3.659"
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -0.97%  '
$ws.Range('D34').Value = "'5.545"
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -3.19%  '
$ws.Range('D35').Value = "'0.2250"
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -2.43%  '
$ws.Range('D36').Value = "'0.06453"
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.93%  '
$ws.Range('D37').Value = "'0.02291"
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -1.16%  '
$ws.Range('D38').Value = "'5.007"
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -3.26%  '
$ws.Range('D39').Value = "'8.436"
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -3.78%  '
$ws.Range('D40').Value = "'1.450"
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +4.61%  '
$ws.Range('D41').Value = "'0.6173"
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -3.29%  '
$ws.Range('E42').Value = '  -4.34%  '
$ws.Range('D43').Value = "'1.176"
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +1.79%  '
$ws.Range('D44').Value = "'1.000"
$ws.Range('D44').Style = "Normal"
$ws.Range('D45').Value = "'13.29"
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -2.16%  '
$ws.Range('E46').Value = '  -0.04%  '
$ws.Range('D47').Value = "'0.5777"
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -3.00%  '
$ws.Range('D48').Value = "'124.06"
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -1.60%  '
$ws.Range('D49').Value = "'1.188"
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +3.46%  '
$ws.Range('E50').Value = '  -3.35%  '
$ws.Range('D51').Value = "'0.06807"
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.29%  '
